$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.684.69"
$ws.Range("D3").Value = "1.808.16"
$ws.Range("E3").Value = "  -1.84%  "
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").Value = "'232.31"
$ws.Range("E5").Value = "  +1.23%  "
$ws.Range("D6").Value = "'0.603"
$ws.Range("E6").Value = "  -1.11%  "
$ws.Range("E7").Value = "  +0.32%  "
$ws.Range("D8").Value = "'39.35"
$ws.Range("E8").Value = "  -9.14%  "
$ws.Range("E9").Value = "  +4.85%  "
$ws.Range("D10").Value = "'0.0681"
$ws.Range("E10").Value = "  -2.01%  "
$ws.Range("E11").Value = "  -1.79%  "
$ws.Range("D12").Value = "2.069.28"
$ws.Range("E12").Value = "  -1.86%  "
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "'0.674"
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.796.27"
$ws.Range("E14").Value = "  -2.49%  "
$ws.Range("D15").Value = "'11.02"
$ws.Range("E15").Value = "  -2.46%  "
$ws.Range("D16").Value = "'4.58"
$ws.Range("E16").Value = "  -1.97%  "
$ws.Range("D17").Value = "34.707.89"
$ws.Range("E17").Value = "  -1.91%  "
$ws.Range("D18").Value = "'69.48"
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("D19").Value = "0.0₃0785"
$ws.Range("E19").Value = "  -1.14%  "
$ws.Range("D20").Value = "'239.87"
$ws.Range("E20").Value = "  -1.87%  "
$ws.Range("D21").Value = "'11.92"
$ws.Range("E21").Value = "  -1.56%  "
$ws.Range("D22").Value = "'4.65"
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("E24").Value = "  +2.20%  "
$ws.Range("D25").Value = "'171.92"
$ws.Range("D26").Value = "'7.71"
$ws.Range("E26").Value = "  -2.36%  "
$ws.Range("D27").Value = "'17.16"
$ws.Range("E27").Value = "  -3.43%  "
$ws.Range("E28").Value = "  -1.83%  "
$ws.Range("D29").Value = "'1.55"
$ws.Range("E29").Value = "  +11.56%  "
$ws.Range("E30").Value = "  +0.41%  "
$ws.Range("D31").Value = "'4.02"
$ws.Range("E31").Value = "  +2.52%  "
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("D33").Value = "'3.97"
$ws.Range("E33").Value = "  -2.50%  "
$ws.Range("D34").Value = "'1.28"
$ws.Range("E34").Value = "  +18.20%  "
$ws.Range("D35").Value = "'1.78"
$ws.Range("E35").Value = "  -4.58%  "
$ws.Range("D36").Value = "'0.700"
$ws.Range("E36").Value = "  +1.67%  "
$ws.Range("D37").Value = "'91.50"
$ws.Range("E37").Value = "  -3.56%  "
$ws.Range("E38").Value = "  +5.64%  "
$ws.Range("D39").Value = "1.319.75"
$ws.Range("E39").Value = "  -1.96%  "
$ws.Range("E40").Value = "  -1.27%  "
$ws.Range("D41").Value = "'2.48"
$ws.Range("E41").Value = "  +0.78%  "
$ws.Range("E42").Value = "  -4.10%  "
$ws.Range("D43").Value = "'14.29"
$ws.Range("E43").Value = "  -3.34%  "
$ws.Range("E45").Value = "  -5.29%  "
$ws.Range("D46").Value = "'6.24"
$ws.Range("E46").Value = "  -0.20%  "
$ws.Range("D47").Value = "'0.0513"
$ws.Range("E47").Value = "  -1.30%  "
$ws.Range("D48").Value = "1.997.18"
$ws.Range("E48").Value = "  -0.60%  "
$ws.Range("E49").Value = "  +0.34%  "
$ws.Range("E50").Value = "  +7.30%  "
$ws.Range("E51").Value = "  -4.39%  "
